$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11..114 down to 12..115.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new data record.
$ws.Range("A11").Value = 3
$ws.Range("B11").Value = "Femacal de La Calera"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44552
$ws.Range("D11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 100112030
$ws.Range("G11").Value = "Poroto granado"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 75
$ws.Range("K11").Value = 44000
$ws.Range("L11").Value = 45000
$ws.Range("M11").Value = 44533
$ws.Range("N11").Value = "$/saco 25 kilos"
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 1781
$ws.Range("Q11").Value = 25
$ws.Range("R11").Value = "Hortaliza"
